$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 337 (shifts rows 337:416 down to 338:417)
$ws.Rows.Item(337).Insert()

$ws.Cells.Item(337, 1).Value = "ritzbet"
$ws.Cells.Item(337, 2).Value = "Maks 500tl çekim"
$ws.Cells.Item(337, 3).Value = "yatırımsız"
$ws.Cells.Item(337, 4).Value = "Evet"

# Match the author's final view state (scroll position / active selection)
$null = $ws.Range("B9").Select()
